# Fruta / hortaliza, semanal
# Insert a new weekly record at row 42 (Macroferia Regional de Talca -
# Arándano (blue)), which pushes all subsequent rows (old 42..57) down by
# one to become rows 43..58.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 42, shifting rows 42:57 down to 43:58.
$ws.Rows("42:42").Insert()

# Populate the newly inserted row 42 with this week's record.
$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Macroferia Regional de Talca"
$ws.Range("C42").Value = "Maule"
$ws.Range("D42").Value = 44588
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100101
$ws.Range("H42").Value = "Berries"
$ws.Range("I42").Value = 100101001
$ws.Range("J42").Value = "Arándano (blue)"
$ws.Range("K42").Value = "Sin especificar"
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 150
$ws.Range("N42").Value = 3500
$ws.Range("O42").Value = 3500
$ws.Range("P42").Value = 3500
$ws.Range("Q42").Value = '$/bandeja 2 kilos'
$ws.Range("R42").Value = "Provincia de Linares"
$ws.Range("S42").Value = 1750
$ws.Range("T42").Value = 2
